$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "whether it be the prerecorded",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "whether it is the prerecorded",
    2
)
